$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so values like "1.001"
# are stored verbatim (inline/shared string) instead of being coerced to a
# number by Excel's smart-entry parsing. Restore style afterwards so no
# extra style index is left on the cells (matches original unstyled cells).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.936.99'
$ws.Range("E2").Value = '  -2.37%  '
$ws.Range("D3").Value = '1.795.24'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '316.87'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.5304'
$ws.Range("E7").Value = '  -4.04%  '
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").Value = '0.07441'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = '41.45'
$ws.Range("E10").Value = '  -2.17%  '
$ws.Range("D11").Value = '1.085'
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '6.167'
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '7.441'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("D16").Value = '1.800.79'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '88.28'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '0.06555'
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '17.21'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = '5.941'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '27.969.44'
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '2.091'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '157.04'
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").Value = '20.09'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").Value = '2.004.56'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '2.292'
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("D30").Value = '121.63'
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.099'
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1085'
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").Value = '3.673'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("D34").Value = '5.481'
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").Value = '0.07090'
$ws.Range("E35").Value = '  +8.82%  '
$ws.Range("D36").Value = '0.2195'
$ws.Range("E36").Value = '  -2.52%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.102'
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02274'
$ws.Range("E38").Value = '  -1.45%  '
$ws.Range("D39").Value = '8.395'
$ws.Range("E39").Value = '  -4.12%  '
$ws.Range("D40").Value = '11.21'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6106'
$ws.Range("E41").Value = '  -2.16%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.179'
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("D43").Value = '1.418'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '3.675'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '0.5697'
$ws.Range("E46").Value = '  -2.71%  '
$ws.Range("D47").Value = '125.03'
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = '1.176'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.913'
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").Value = '0.06805'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = '0.00000000293'
$ws.Range("E51").Value = '  +25.21%  '

$dataRange.Style = "Normal"
